$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.981.45'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '1.593.86'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  -0.93%  '
$ws.Range("E9").Value = '  -1.23%  '
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("E11").Value = '  +2.79%  '
$ws.Range("D12").Value = '1.815.85'
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").Value = '1.589.48'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").Value = '25.984.96'
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '199.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.122'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.07%  '
$ws.Range("E33").Value = '  -3.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("D36").Value = '1.123.38'
$ws.Range("E36").Value = '  +1.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0163'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.45%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -1.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.783'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.490'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.783'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.20%  '
$ws.Range("D43").Value = '1.726.86'
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '92.56'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.09'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.27%  '
$ws.Range("E46").Value = '  -0.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '53.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("D51").Value = '0.0₇0917'
$ws.Range("E51").Value = '  -18.15%  '
